$d = $word.ActiveDocument

# 1) "...EL CUAL TIENEN LAS SIGUIENTES MEDIDA Y COLINDANCIAS" -> "...MEDIDAS Y COLINDANCIAS"
$d.Content.Find.Execute("SIGUIENTES MEDIDA Y COLINDANCIAS", $true, $false, $false, $false, $false, `
    $true, 1, $false, "SIGUIENTES MEDIDAS Y COLINDANCIAS", 2) | Out-Null

# 2) Drop the opening curly left-quote before {{SEXO_7}} and collapse the
#    space between {{SEXO_7}} and PROMITENTE (8 occurrences across the doc;
#    the matching closing curly right-quote is left untouched).
$leftQuote = [char]0x201C
$findSexo7 = $leftQuote + "{{SEXO_7}} PROMITENTE"
$replaceSexo7 = "{{SEXO_7}}PROMITENTE"
$d.Content.Find.Execute($findSexo7, $true, $false, $false, $false, $false, `
    $true, 1, $false, $replaceSexo7, 2) | Out-Null

# 3) {{SEXO_17}} -> ÉSTA ; {{SEXO_16}} -> ÉSTA
$d.Content.Find.Execute("{{SEXO_17}}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ÉSTA", 2) | Out-Null
$d.Content.Find.Execute("{{SEXO_16}}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ÉSTA", 2) | Out-Null
